$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. optimization_parameters sheet: restructure rows 1 and 8-16
# ---------------------------------------------------------------------------
$wsOpt = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the redundant C1:F1 cells (they duplicated B1's "value" label)
$wsOpt.Range("C1:F1").ClearContents()

# Insert a new row below the "Model" row (row 8) to hold the new "L_curve" parameter
$wsOpt.Rows.Item(9).Insert()

# Rename the "Model" label (row 8) to "production_function" - value stays "Sigmoid"
$wsOpt.Range("A8").Value = "production_function"

# Populate the newly inserted row 9 with the "L_curve" parameter
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 1
$wsOpt.Range("B9").NumberFormat = "0.00E+00"

# Remove the obsolete "Deletion" row (now at row 17, after the insert above)
$wsOpt.Rows.Item(17).Delete()

# ---------------------------------------------------------------------------
# 2. Switch the active sheet / selection from dcin5_log2_expression to
#    optimization_parameters
# ---------------------------------------------------------------------------
$wsOpt.Activate()
$wsOpt.Range("C1:F4").Select()
